$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, pushing existing rows 3-9 down to 4-10
$ws.Rows.Item(3).Insert()

# Fill in the new row 3 with the new weekly data point.
# Columns A,B,C,E,F,G,H,I,J are constant across all rows of this sheet.
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = 44544
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100103
$ws.Cells.Item(3, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(3, 9).Value = 100103003
$ws.Cells.Item(3, 10).Value = "Damasco"
$ws.Cells.Item(3, 11).Value = "Castle Brite"
$ws.Cells.Item(3, 12).Value = "Segunda"
$ws.Cells.Item(3, 13).Value = 160
$ws.Cells.Item(3, 14).Value = 16000
$ws.Cells.Item(3, 15).Value = 17000
$ws.Cells.Item(3, 16).Value = 16500
$ws.Cells.Item(3, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(3, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(3, 19).Value = 1100
$ws.Cells.Item(3, 20).Value = 15
